# The workbook is a small 3-column table (收款方/付款方/金额) backed by a
# shared-strings table. The edit:
#   - changes cell B2's text from "aaa" to a new string "fff"
#   - leaves the active selection on B2 (matching the saved cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "fff"

$ws.Range("B2").Select()
